# 17.1.1 - add year 2023 (column T) and update 2022 (column S) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- column widths (engine rounds ColumnWidth to 1/6-character increments, so
#    these land as close as the host allows to the authored 8.5703125 width) --
$ws.Range("D1:L1").ColumnWidth = 7.6
$ws.Range("M1:O1").ColumnWidth = 7.6
$ws.Range("P1:T1").ColumnWidth = 7.6

# -- row heights --
$ws.Rows.Item(1).RowHeight = 42.75

# -- existing column S (2022) figure corrections --
$ws.Range("S5").Value = 29.5
$ws.Range("S6").Value = 22.4
$ws.Range("S8").Value = 1.9
$ws.Range("S9").Value = 5.0999999999999996

# -- new column T (2023) --

# Row 4: header year
$ws.Range("T4").Value = 2023
$ws.Range("T4").Font.Bold = $true
$ws.Range("T4").Font.Size = 9
$ws.Range("T4").Font.Name = "Times New Roman"
$ws.Range("T4").HorizontalAlignment = -4152
$ws.Range("T4").VerticalAlignment = -4108
$ws.Range("T4").Borders.Item(8).LineStyle = 1
$ws.Range("T4").Borders.Item(8).Weight = -4138
$ws.Range("T4").Borders.Item(9).LineStyle = 1
$ws.Range("T4").Borders.Item(9).Weight = -4138

# Row 5: Revenues, total
$ws.Range("T5").Value = 29.4
$ws.Range("T5").Font.Bold = $true
$ws.Range("T5").Font.Size = 9
$ws.Range("T5").Font.Name = "Times New Roman"
$ws.Range("T5").NumberFormat = "0.0"
$ws.Range("T5").HorizontalAlignment = -4152
$ws.Range("T5").VerticalAlignment = -4108

# Row 6: Tax revenues
$ws.Range("T6").Value = 22.1
$ws.Range("T6").Font.Bold = $false
$ws.Range("T6").Font.Size = 9
$ws.Range("T6").Font.Name = "Times New Roman"
$ws.Range("T6").NumberFormat = "0.0"
$ws.Range("T6").HorizontalAlignment = -4152
$ws.Range("T6").VerticalAlignment = -4108

# Row 7: Contributions / deductions for social needs
$ws.Range("T7").Value = "-"
$ws.Range("T7").Font.Bold = $false
$ws.Range("T7").Font.Size = 9
$ws.Range("T7").Font.Name = "Times New Roman"
$ws.Range("T7").HorizontalAlignment = -4152
$ws.Range("T7").VerticalAlignment = -4108

# Row 8: Received official transfers
$ws.Range("T8").Value = 1.2
$ws.Range("T8").Font.Bold = $false
$ws.Range("T8").Font.Size = 9
$ws.Range("T8").Font.Name = "Times New Roman"
$ws.Range("T8").NumberFormat = "0.0"
$ws.Range("T8").HorizontalAlignment = -4152
$ws.Range("T8").VerticalAlignment = -4108

# Row 9: Non-tax revenues
$ws.Range("T9").Value = 6.1
$ws.Range("T9").Font.Bold = $false
$ws.Range("T9").Font.Size = 9
$ws.Range("T9").Font.Name = "Times New Roman"
$ws.Range("T9").NumberFormat = "0.0"
$ws.Range("T9").HorizontalAlignment = -4152
$ws.Range("T9").VerticalAlignment = -4108

# Row 10: Revenues from the sale of non-financial assets
$ws.Range("T10").Value = 0
$ws.Range("T10").Font.Bold = $false
$ws.Range("T10").Font.Size = 9
$ws.Range("T10").Font.Name = "Times New Roman"
$ws.Range("T10").NumberFormat = "0.0"
$ws.Range("T10").HorizontalAlignment = -4152
$ws.Range("T10").VerticalAlignment = -4108
$ws.Range("T10").Borders.Item(9).LineStyle = 1
$ws.Range("T10").Borders.Item(9).Weight = -4138
